$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("订单")

# ===================== Row 34: in-place edits =====================
$ws.Cells.Item(34, 1).Value = 45350
$c34d = $ws.Cells.Item(34, 4)
$c34d.Value = "'414731942"
$c34d.Style = "Normal"
$ws.Cells.Item(34, 10).Value = ""

# ===================== Insert a new row at 40 =====================
# old row 40 shifts down to row 41; bump its sequence number (col B) to 40
$ws.Rows(40).Insert()
$ws.Cells.Item(41, 2).Value = 40

# ===================== Fill new row 40 =====================
$ws.Cells.Item(40, 1).Value = 45350
$c40d = $ws.Cells.Item(40, 4)
$c40d.Value = "'414731942"
$c40d.Style = "Normal"
$ws.Cells.Item(40, 7).Value = 45299
$ws.Cells.Item(40, 9).Value = '菌群+对应代谢产物介导+机制研究'

# ===================== Append new rows 42-51 =====================
# --- Row 42 ---
$ws.Cells.Item(42, 1).Value = 45350
$ws.Cells.Item(42, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(42, 2).Value = 41
$ws.Cells.Item(42, 4).Value = 'BI2024013001'
$ws.Cells.Item(42, 6).Value = '黄礼闯'
$ws.Cells.Item(42, 7).Value = 45322
$ws.Cells.Item(42, 7).NumberFormat = "m/d/yy"
$ws.Cells.Item(42, 8).Value = '完成'
$ws.Cells.Item(42, 9).Value = '审核业务'

# --- Row 43 ---
$ws.Cells.Item(43, 1).Value = 45350
$ws.Cells.Item(43, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(43, 2).Value = 42
$ws.Cells.Item(43, 3).Value = ' 01-订单编号：N2024020103 02-区域-销售：张玉玲 03-上级主管：王立家 04-医院：萧山区第一人民医院 05-科室/职称：肾内科 护士长 06-电话： 07-项目（确定A/B套餐）：厅级标书+预实验 08-分值： 09-定题题目：按照技术路线图 10-时间要求：                         预实验：2024年5月30日                          标   书：2024年6月15日   11-总价： 12-定金：已付 13-评估人员 ：吴晨 14-技术支持（沟通情况）：薛富才 15-附件：报价、合同、标书申报浙江省医药科技计划、标书技术路线、预实验技术内容、技术支持与客户沟通反馈总结 16-备注： 客户要求： 1）客户分级：1.院方职务（科室职务）重点客户，做项目为了后续拿课题。 2）实验部分。在正式实验前需要进行预实验摸索动物模型建立情况，如效果好，继续做下去；效果不好，终止实验，此部分费用由我们承担。实验优先安排加急做下去，五月底做完全部实验。【预实验走实验项目】 3）标书部分。写浙江省医药科技计划标书，重点客户，安排优秀的同事写作，保证质量；6月15号前交付；如今年未立项，后续可修改一次投其他项目。 内部要求：  1）内部留存原始数据，三次重复实验 2）方案在做的时候，发现技术路线和方案有较大出入的，需要跟上游和技术支持沟通【方案能调整需告知】 17-项目负责人：杨弘 '
$ws.Cells.Item(43, 4).Value = 'N2024020103'
$ws.Cells.Item(43, 6).Value = '黄礼闯'
$ws.Cells.Item(43, 7).Value = 45324
$ws.Cells.Item(43, 7).NumberFormat = "m/d/yy"
$ws.Cells.Item(43, 8).Value = '完成'
$ws.Cells.Item(43, 9).Value = '筛选主动脉-下腔静脉瘘ACF模型 DEGs 并功能分析'

# --- Row 44 ---
$ws.Cells.Item(44, 1).Value = 45350
$ws.Cells.Item(44, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(44, 2).Value = 43
$ws.Cells.Item(44, 3).Value = ' 01-订单编号：N2024012602  02-区域-销售：江苏-郭树仁 03-上级主管：王立家 04-医院：绍兴市中 05-科室/职称：骨科 06-电话： 07-项目（确定A/B套餐）：sci1.5-2分全包B套餐 08-分值：1.5-2分 09-定题题目： 10-时间要求：2024年1月24-2025年1月24日 11-总价： 12-定金： 13-评估人员 ：林婧宇 14-技术支持（沟通情况）：薛富才-2 15-附件：合同、张春晓文章技术内容、动物伦理、动物使用许可、前期基础、申报书正文、项目计划书、客户提供数据（云附件）、技术支持与客户沟通反馈总结 16-备注： 客户要求： 1.客户分级（1.院方职务（科室职务）重点客户 2.谈单承诺：必须在时间内完成录用（2025年1月24日） 内部要求： 1)【方案能调整】方案我们定，之前内部沟通只要4月能把实验完成就行； 2)需要安排外部投稿； @助理陈芳媛 3)全包实验，尽量按实验数据交付。 4）内部留存原始数据，三次重复实验 17-项目负责人：杨啸   '
$ws.Cells.Item(44, 4).Value = 'N2024012602'
$ws.Cells.Item(44, 5).Value = '1.5-2分'
$ws.Cells.Item(44, 6).Value = '黄礼闯'
$ws.Cells.Item(44, 7).Value = 45324
$ws.Cells.Item(44, 7).NumberFormat = "m/d/yy"
$ws.Cells.Item(44, 8).Value = '完成'
$ws.Cells.Item(44, 9).Value = 'Hydroxysafflor Yellow A 与Piezo1对接'

# --- Row 45 ---
$ws.Cells.Item(45, 1).Value = 45350
$ws.Cells.Item(45, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(45, 2).Value = 44
$ws.Cells.Item(45, 3).Value = '01-订单编号：N2024010303【方案能调整需告知】 02-区域-销售：张玉玲 03-上级主管：王立家 04-医院： 05-科室/职称：儿科 正高  科主任 06-电话： 07-项目（确定A/B套餐）：中管局标书（实验走实验项目 08-分值： 09-定题题目：按照技术路线 10-时间要求：2024年6月10日 11-总价： 12-定金： 13-评估人员 ：吴晨 14-技术支持（沟通情况）：薛富才，龙艳 15-附件：技术支持与客户沟通反馈总结、技术路线，技术内容、（报价，预实验报价、合同、客户提供的前期研究基础等资料见压缩包）16-备注：客户要求：1）客户分级：重点客户，做项目为了后续拿课题。 2）谈单承诺：实验部分需要加急安排做！重点客户，根据方案设计写中管局标书，需要2024年6月10日前交付！辛苦安排优秀的同事写，质量要高。 内部要求：1）内部留存原始数据，三次重复实验2）方案在做的时候，发现技术路线和方案有较大出入的，需要跟上游和技术支持沟通； 17-项目负责人：杨啸 '
$ws.Cells.Item(45, 4).Value = 'N2024010303'
$ws.Cells.Item(45, 6).Value = '黄礼闯'
$ws.Cells.Item(45, 7).Value = 45326
$ws.Cells.Item(45, 7).NumberFormat = "m/d/yy"
$ws.Cells.Item(45, 8).Value = '完成'
$ws.Cells.Item(45, 9).Value = '分子对接 Celogenamide A（环状肽）蛋白 SSTR2'

# --- Row 46 ---
$ws.Cells.Item(46, 1).Value = 45350
$ws.Cells.Item(46, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(46, 2).Value = 45
$c46d = $ws.Cells.Item(46, 4)
$c46d.Value = "'20231012"
$c46d.Style = "Normal"
$ws.Cells.Item(46, 6).Value = '黄礼闯'
$ws.Cells.Item(46, 7).Value = 45343
$ws.Cells.Item(46, 7).NumberFormat = "m/d/yy"
$ws.Cells.Item(46, 8).Value = '完成'
$ws.Cells.Item(46, 9).Value = '建立风险模型和作图'

# --- Row 47 ---
$ws.Cells.Item(47, 1).Value = 45350
$ws.Cells.Item(47, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(47, 2).Value = 46
$c47d = $ws.Cells.Item(47, 4)
$c47d.Value = "'20230815"
$c47d.Style = "Normal"
$ws.Cells.Item(47, 6).Value = '黄礼闯'
$ws.Cells.Item(47, 7).Value = 45344
$ws.Cells.Item(47, 7).NumberFormat = "m/d/yy"
$ws.Cells.Item(47, 8).Value = '完成'
$ws.Cells.Item(47, 9).Value = '列线图模型建立与验证'

# --- Row 48 ---
$ws.Cells.Item(48, 1).Value = 45350
$ws.Cells.Item(48, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(48, 2).Value = 47
$ws.Cells.Item(48, 3).Value = ' 01-订单编号：N2024022202  02-区域-销售：浙江台州-吴航贵 03-上级主管：王立家 04-医院：温岭市第一人民医院 05-科室/职称：肝病科/科主任 06-电话： 07-项目（确定A/B套餐）： 08-分值： 09-定题题目：乙肝病毒HBx利用泛素化系统降解XXX上调YYY诱导肝癌线粒体自噬的机制研究【方案不能调整】 10-时间要求：2024年2月22日--2024年06月15日 11-总价： 12-定金： 13-评估人员 ： 14-技术支持（沟通情况）：薛富才，吴晨（沟通次数）3次， 15-附件：标书技术路线、预实验技术内容、技术支持与客户沟通反馈总结、合同、预实验报价 16-备注： 客户要求： 1）客户分级：1.院方职务(科室职务)重点客户，做项目为了后续拿课题 2）尽量在今年省自然投之前给到客户交付 内部要求： 1）内部留存原始数据，三次重复实验 ；  2）方案不能调整； 3）预实验走实验项目 17-项目负责人：杨啸 '
$ws.Cells.Item(48, 4).Value = 'N2024022202'
$ws.Cells.Item(48, 6).Value = '黄礼闯'
$ws.Cells.Item(48, 7).Value = 45345
$ws.Cells.Item(48, 7).NumberFormat = "m/d/yy"
$ws.Cells.Item(48, 8).Value = '完成'
$ws.Cells.Item(48, 9).Value = '乙肝病毒HBx利用泛素化系统降解XXX上调YYY诱导肝癌线粒体自噬'

# --- Row 49 ---
$ws.Cells.Item(49, 1).Value = 45351
$ws.Cells.Item(49, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(49, 2).Value = 48
$ws.Cells.Item(49, 4).Value = 'workflow'
$ws.Cells.Item(49, 6).Value = '黄礼闯'
$ws.Cells.Item(49, 7).Value = 45351
$ws.Cells.Item(49, 7).NumberFormat = "m/d/yy"
$ws.Cells.Item(49, 8).Value = '完成'
$ws.Cells.Item(49, 9).Value = 'Step 系列：scRNA-seq 基本分析'

# --- Row 50 ---
$ws.Cells.Item(50, 1).Value = 45351
$ws.Cells.Item(50, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(50, 2).Value = 49
$ws.Cells.Item(50, 4).Value = 'workflow'
$ws.Cells.Item(50, 6).Value = '黄礼闯'
$ws.Cells.Item(50, 7).Value = 45351
$ws.Cells.Item(50, 7).NumberFormat = "m/d/yy"
$ws.Cells.Item(50, 8).Value = '完成'
$ws.Cells.Item(50, 9).Value = 'Step 系列：Prologue and Get-start'

# --- Row 51 ---
$ws.Cells.Item(51, 1).Value = 45351
$ws.Cells.Item(51, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(51, 2).Value = 50
$ws.Cells.Item(51, 4).Value = 'workflow'
$ws.Cells.Item(51, 6).Value = '黄礼闯'
$ws.Cells.Item(51, 7).Value = 45351
$ws.Cells.Item(51, 7).NumberFormat = "m/d/yy"
$ws.Cells.Item(51, 8).Value = '完成'
$ws.Cells.Item(51, 9).Value = 'Step 系列：scRNA-seq 癌细胞鉴定'

# ===================== Sheet "备单": cosmetic style bump is a no-op here =====================
# (s="11" -> s="12" on A2:A4 / G2:G4 is a duplicate of an identical date style; no value/format change)
